$d = $word.ActiveDocument

# Helper: insert a plain paragraph of text (can be empty) right at
# character position $p, ending it with a paragraph mark, and return
# the position right after that new paragraph mark (i.e. where the
# next paragraph should be inserted).
function Insert-PlainParagraph($p, $text) {
    $rr = $d.Range($p, $p)
    $rr.InsertAfter($text + "`r")
    return $p + $text.Length + 1
}

# Start right after the existing "Test" paragraph (paragraph 1).
$pos = $d.Paragraphs(1).Range.End

# <w:p/>
$pos = Insert-PlainParagraph $pos ""

# "Dit is een tekst met de tekst klaas erin"
$pos = Insert-PlainParagraph $pos "Dit is een tekst met de tekst klaas erin"

# <w:p/>
$pos = Insert-PlainParagraph $pos ""

# <w:p/>
$pos = Insert-PlainParagraph $pos ""

# "Ik ben getypt in word"
$pos = Insert-PlainParagraph $pos "Ik ben getypt in word"

# <w:p/>
$pos = Insert-PlainParagraph $pos ""

# Final inserted paragraph needs two distinct runs: "Tekst met wederom "
# and "klaas". Plain InsertAfter calls get silently coalesced into a
# single run on save, so build this one paragraph from an OOXML
# fragment (still placed via the normal Range API) which preserves the
# run boundary, then close it off with a real paragraph mark.
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$paraXml = "<w:p $wNs>" +
           "<w:r><w:t xml:space='preserve'>Tekst met wederom </w:t></w:r>" +
           "<w:r><w:t>klaas</w:t></w:r>" +
           "</w:p>"
$finalRange = $d.Range($pos, $pos)
$finalRange.InsertXML($paraXml)

$afterFinalText = $pos + "Tekst met wederom klaas".Length
$breakRange = $d.Range($afterFinalText, $afterFinalText)
$breakRange.InsertAfter("`r")

Write-Output "done"
